# Apply the cryptos.xlsx data refresh described in the commit diff.
# Column D ("Price") cells sometimes look numeric (e.g. "1.00", "0.999",
# "0.0000216"); Excel would silently coerce those to real numbers on a
# plain .Value assignment, which changes their stored type/formatting away
# from the original plain-text cells. Prefixing with a leading apostrophe
# forces text entry, and resetting .Style back to "Normal" afterwards
# strips the text-number-format style Excel auto-applies, so the cell ends
# up as plain text with the default style - matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''61.493.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.67%  '

# Row 3
$ws.Range('D3').Value = '''2.972.63'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.06%  '

# Row 4
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.20%  '

# Row 5
$ws.Range('D5').Value = '''527.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.41%  '

# Row 6
$ws.Range('D6').Value = '''129.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.28%  '

# Row 7
$ws.Range('D7').Value = '''1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('D8').Value = '''2.968.74'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.15%  '

# Row 9
$ws.Range('D9').Value = '''0.485'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.45%  '

# Row 10
$ws.Range('E10').Value = '  -2.91%  '

# Row 11
$ws.Range('D11').Value = '''6.08'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.38%  '

# Row 12
$ws.Range('D12').Value = '''0.438'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.71%  '

# Row 13
$ws.Range('D13').Value = '''0.0000216'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.72%  '

# Row 14
$ws.Range('D14').Value = '''33.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.42%  '

# Row 15
$ws.Range('D15').Value = '''3.477.09'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.64%  '

# Row 16
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value = '''0.110'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.03%  '

# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '''61.455.36'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.99%  '

# Row 18
$ws.Range('D18').Value = '''2.975.31'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.25%  '

# Row 19
$ws.Range('D19').Value = '''6.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.86%  '

# Row 20
$ws.Range('D20').Value = '''454.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.56%  '

# Row 21
$ws.Range('D21').Value = '''13.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.56%  '

# Row 22
$ws.Range('D22').Value = '''0.667'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.95%  '

# Row 23
$ws.Range('D23').Value = '''6.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.69%  '

# Row 24
$ws.Range('D24').Value = '''77.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.87%  '

# Row 25
$ws.Range('D25').Value = '''11.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.85%  '

# Row 26
$ws.Range('D26').Value = '''0.996'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.46%  '

# Row 27
$ws.Range('E27').Value = '  -1.88%  '

# Row 28
$ws.Range('D28').Value = '''7.55'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.84%  '

# Row 29
$ws.Range('D29').Value = '''1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.08%  '

# Row 30
$ws.Range('D30').Value = '''25.26'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.43%  '

# Row 31
$ws.Range('E31').Value = '  +2.13%  '

# Row 32
$ws.Range('D32').Value = '''1.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.07%  '

# Row 33
$ws.Range('D33').Value = '''55.65'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.01%  '

# Row 34
$ws.Range('E34').Value = '  -6.90%  '

# Row 35
$ws.Range('D35').Value = '''5.27'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.82%  '

# Row 36
$ws.Range('D36').Value = '''5.73'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.91%  '

# Row 37
$ws.Range('D37').Value = '''450.10'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.39%  '

# Row 38
$ws.Range('D38').Value = '''3.121.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.57%  '

# Row 39
$ws.Range('D39').Value = '''0.0381'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.13%  '

# Row 40
$ws.Range('D40').Value = '''0.0772'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.82%  '

# Row 41
$ws.Range('D41').Value = '''0.114'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.87%  '

# Row 42
$ws.Range('D42').Value = '''7.90'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.24%  '

# Row 43
$ws.Range('D43').Value = '''2.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.40%  '

# Row 45
$ws.Range('E45').Value = '  -2.88%  '

# Row 46
$ws.Range('D46').Value = '''24.67'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.70%  '

# Row 47
$ws.Range('D47').Value = '''119.61'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.44%  '

# Row 48
$ws.Range('E48').Value = '  -0.33%  '

# Row 49
$ws.Range('B49').Value = 'PEPE'
$ws.Range('C49').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D49').Value = '''0.0₃0498'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.84%  '

# Row 50
$ws.Range('B50').Value = 'Fetch.AI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D50').Value = '''1.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.78%  '

# Row 51
$ws.Range('D51').Value = '''1.23'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.48%  '
